# Scene 20Aa edit: merge the split runs in the "Asher didn't show up..."
# paragraph back into a single contiguous run of text (no wording change,
# just collapsing runs that had been split across several <w:r> elements).

$d = $word.ActiveDocument

$oldText = "Asher didn’t show up today, so after barely making it through class I make my way to the stairwell at the end of school, a small part of me hoping that Lilith will be there. However, she isn’t, and in embarrassment I quickly leave, resolving to eat in the classroom."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

if (-not $result) {
    throw "Could not find the target paragraph text to normalize."
}
